$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51..171 down to 52..172
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new weekly price record
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44581
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100112005
$ws.Range("G51").Value = "Puerro"
$ws.Range("H51").Value = "Azul de Maquehue"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 100
$ws.Range("K51").Value = 12000
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = 12000
$ws.Range("N51").Value = "`$/docena de paquetes"
$ws.Range("O51").Value = "Provincia de Cautín"
$ws.Range("P51").Value = 1000
$ws.Range("Q51").Value = 12
$ws.Range("R51").Value = "Hortaliza"
